$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17: clear A17's value (keep style) and remove B17 / C17 entirely
$ws.Range("A17").ClearContents()
$ws.Range("B17:C17").Clear()

# Row 18: new values, drop the shared-formula C18
$ws.Range("A18").Value = 4286
$ws.Range("B18").Value = 207
$ws.Range("C18").Clear()

# Row 19: new values, drop the shared-formula C19
$ws.Range("A19").Value = 3960
$ws.Range("B19").Value = 191
$ws.Range("C19").Clear()

# Row 20: previously empty, now gets values (B20 needs the same style as B18/B19, s="1")
$ws.Range("A20").Value = 4000
$ws.Range("B20").Value = 195
$ws.Range("B20").Font.Size = 12

# Update the selected cell shown when the workbook is opened
$ws.Range("B23").Select()
